$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 263, pushing the existing row 263..297 down to 264..298
$ws.Rows(263).Insert()

# Populate the new row 263 with the new price-report record
$ws.Range("A263").Value = 9
$ws.Range("B263").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C263").Value = "Metropolitana"
$ws.Range("D263").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D263").Value = 44476
$ws.Range("E263").Value = 13
$ws.Range("F263").Value = 100112024
$ws.Range("G263").Value = "Choclo"
$ws.Range("H263").Value = "Dulce o Americano"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 52
$ws.Range("K263").Value = 43000
$ws.Range("L263").Value = 44000
$ws.Range("M263").Value = 43500
$ws.Range("N263").Value = "`$/malla 70 unidades"
$ws.Range("O263").Value = "Región de Arica y Parinacota"
$ws.Range("P263").Value = 621
$ws.Range("Q263").Value = 70
$ws.Range("R263").Value = "Hortaliza"
